# feat: add 2022-Q4 data
#
# Workbook originally has 2 sheets: "总计" (summary) and "2021-Q4" (fund
# holdings detail for 2021-Q4). This adds a new "2022-Q4" quarter:
#   - "总计" gets a new top row for 2022-Q4 (pushing the 2021-Q4 row down)
#   - the existing "2021-Q4" detail sheet (currently sheetId=2 / rId2)
#     is repurposed to hold the NEW 2022-Q4 detail data and renamed
#     "2022-Q4" (so it keeps sheetId=2, matching the target layout)
#   - a brand new sheet is appended at the end, named "2021-Q4", holding
#     the OLD fund-holdings data that used to live in the detail sheet
#     (this new sheet becomes sheetId=3 / rId3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert the 2022-Q4 row above the existing 2021-Q4 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy row 2's formatting down into row 3 first, so the new row3 carries
# the same styling (border/bold/center) that row 2 has on column A.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.37

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 1.21

# ---------------------------------------------------------------------
# 2) Repurpose the existing "2021-Q4" detail sheet -> "2022-Q4"
#    (keeps its current sheetId/rId, matching target's rId2/sheetId2)
# ---------------------------------------------------------------------
$detail2022 = $wb.Worksheets.Item(2)
$detail2022.Cells.Clear()
$detail2022.Name = "2022-Q4"

$detail2022.Range("B1").Value = "基金代码"
$detail2022.Range("C1").Value = "基金名称"
$detail2022.Range("D1").Value = "基金规模"
$detail2022.Range("E1").Value = "股票总仓位"
$detail2022.Range("F1").Value = "仓位占比"
$detail2022.Range("G1").Value = "持有市值(亿元)"
$detail2022.Range("H1").Value = "仓位排名"
$detail2022.Range("B1:H1").Style = "Normal"

$detail2022.Range("A2").Value = 0
$detail2022.Range("B2").Value = "'501087"
$detail2022.Range("C2").Value = "交银施罗德瑞丰混合（LOF）"
$detail2022.Range("D2").Value = "'21.14"
$detail2022.Range("E2").Value = "'88.79"
$detail2022.Range("F2").Value = "'4.67"
$detail2022.Range("G2").Value = "'0.9872"
$detail2022.Range("H2").Value = 9
$detail2022.Range("A2:H2").Style = "Normal"

$detail2022.Range("A3").Value = 1
$detail2022.Range("B3").Value = "'011924"
$detail2022.Range("C3").Value = "嘉实港股互联网产业核心资产混合A"
$detail2022.Range("D3").Value = "'1.79"
$detail2022.Range("E3").Value = "'90.30"
$detail2022.Range("F3").Value = "'7.62"
$detail2022.Range("G3").Value = "'0.1364"
$detail2022.Range("H3").Value = 6
$detail2022.Range("A3:H3").Style = "Normal"

$detail2022.Range("A4").Value = 2
$detail2022.Range("B4").Value = "'011925"
$detail2022.Range("C4").Value = "嘉实港股互联网产业核心资产混合C"
$detail2022.Range("D4").Value = "'1.14"
$detail2022.Range("E4").Value = "'90.30"
$detail2022.Range("F4").Value = "'7.62"
$detail2022.Range("G4").Value = "'0.0869"
$detail2022.Range("H4").Value = 6
$detail2022.Range("A4:H4").Style = "Normal"

# Re-apply the header/index styling (bold, centered, bordered) that the
# original sheet used for row 1 and column A, same style class as the
# "总计" sheet's header/index cells.
$summary.Range("B1").Copy()
$detail2022.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("A2").Copy()
$detail2022.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 3) New sheet at the end holding the OLD 2021-Q4 detail data
# ---------------------------------------------------------------------
$detail2021 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $detail2022)
$detail2021.Name = "2021-Q4"

$detail2021.Range("B1").Value = "基金代码"
$detail2021.Range("C1").Value = "基金名称"
$detail2021.Range("D1").Value = "基金规模"
$detail2021.Range("E1").Value = "股票总仓位"
$detail2021.Range("F1").Value = "仓位占比"
$detail2021.Range("G1").Value = "持有市值(亿元)"
$detail2021.Range("H1").Value = "仓位排名"
$detail2021.Range("B1:H1").Style = "Normal"

$detail2021.Range("A2").Value = 0
$detail2021.Range("B2").Value = "'011153"
$detail2021.Range("C2").Value = "华宝新兴消费混合A"
$detail2021.Range("D2").Value = "'12.59"
$detail2021.Range("E2").Value = "'74.02"
$detail2021.Range("F2").Value = "'2.91"
$detail2021.Range("G2").Value = "'0.3664"
$detail2021.Range("H2").Value = 8
$detail2021.Range("A2:H2").Style = "Normal"

$detail2021.Range("A3").Value = 1
$detail2021.Range("B3").Value = "'011154"
$detail2021.Range("C3").Value = "华宝新兴消费混合C"
$detail2021.Range("D3").Value = "'0.16"
$detail2021.Range("E3").Value = "'74.02"
$detail2021.Range("F3").Value = "'2.91"
$detail2021.Range("G3").Value = "'0.0047"
$detail2021.Range("H3").Value = 8
$detail2021.Range("A3:H3").Style = "Normal"

$summary.Range("B1").Copy()
$detail2021.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("A2").Copy()
$detail2021.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats

$null = $summary.Select()
$null = $summary.Range("A1").Select()
